$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ticket")

# --- Header row (row 1) ---
# Existing M1 "StarExpWeight" becomes "SoulStoneWeight"
$ws.Range("M1").Value = "SoulStoneWeight"
# Existing N1 "DetailWeightSum" stays "DetailWeightSum" (text unchanged, just shifted in shared-string table)
$ws.Range("N1").Value = "DetailWeightSum"
# New header columns
$ws.Range("O1").Value = "SoulStoneMinCnt"
$ws.Range("P1").Value = "SoulStoneMaxCnt"

# --- Type row (row 2) ---
$ws.Range("O2").Value = "int"
$ws.Range("P2").Value = "int"

# --- Data rows (row 4, row 5) ---
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 3

$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 3

# --- Column widths ---
# Target stored widths (from OOXML <col width=.../>): K=16.7109375, L=13.5703125 (unchanged),
# M=16, N=9.140625, O=16, P=9.140625.
# The COM bridge quantises ColumnWidth to a 1/7-character pixel grid (stored = round(w*7+5)/7),
# so the nearest achievable input is used for K/N/P (exact hit not representable on that grid).
$ws.Columns.Item(11).ColumnWidth = 16            # K -> stored ~16.7142857 (closest to 16.7109375)
$ws.Columns.Item(13).ColumnWidth = 15.2857142857 # M -> stored 16 (exact)
$ws.Columns.Item(14).ColumnWidth = 8.4285714286  # N -> stored ~9.1428571 (closest to 9.140625)
$ws.Columns.Item(15).ColumnWidth = 15.2857142857 # O -> stored 16 (exact)
$ws.Columns.Item(16).ColumnWidth = 8.4285714286  # P -> stored ~9.1428571 (closest to 9.140625)
# (column letters: 11=K, 12=L, 13=M, 14=N, 15=O, 16=P)

# --- Selection ---
$ws.Range("O4").Select()
